$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in column F, matching the style used by the other header cells
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Timestamps for rows 2-31 (time_taken metadata values)
$timestamps = @(
    "2021-10-05 13:39:54.269333",
    "2021-10-05 13:39:54.269345",
    "2021-10-05 13:39:54.269349",
    "2021-10-05 13:39:54.269352",
    "2021-10-05 13:39:54.269356",
    "2021-10-05 13:39:54.269359",
    "2021-10-05 13:39:54.269362",
    "2021-10-05 13:39:54.269365",
    "2021-10-05 13:39:54.269369",
    "2021-10-05 13:39:54.269372",
    "2021-10-05 13:39:54.269375",
    "2021-10-05 13:39:54.269378",
    "2021-10-05 13:39:54.269381",
    "2021-10-05 13:39:54.269384",
    "2021-10-05 13:39:54.269387",
    "2021-10-05 13:39:54.269390",
    "2021-10-05 13:39:54.269393",
    "2021-10-05 13:39:54.269396",
    "2021-10-05 13:39:54.269400",
    "2021-10-05 13:39:54.269403",
    "2021-10-05 13:39:54.269406",
    "2021-10-05 13:39:54.269409",
    "2021-10-05 13:39:54.269412",
    "2021-10-05 13:39:54.269415",
    "2021-10-05 13:39:54.269419",
    "2021-10-05 13:39:54.269422",
    "2021-10-05 13:39:54.269425",
    "2021-10-05 13:39:54.269428",
    "2021-10-05 13:39:54.269431",
    "2021-10-05 13:39:54.269434"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
